# HTML - UI Base Design Prepared
# Clear the now-unused "Yes"/"No" values from the Execute_Flag helper column
# (column C) on the TC_Details sheet. This drops the two stray shared
# strings ("Yes"/"No") and shrinks the sheet's used range back down to
# column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC_Details")

$ws.Range("C2").Value = ""
$ws.Range("C3").Value = ""
